# Update "Forecast Comparison" sheet with corrected forecast output:
#  - Insert a new "Week_Start_Date" column after "Week" (new column B)
#  - Shorten week labels from "W01".."W09" to "W1".."W9" (W10+ stay as-is)
#  - Populate the new Week_Start_Date column with the week's start date (as text)
#  - Convert is_holiday_week column (now column J) values to real booleans

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1. Insert a new column before the current column B (ASIN), shifting
#    ASIN..is_holiday_week one column to the right (B->C, C->D, ..., I->J)
$ws.Columns("B:B").Insert()

# 2. Header for the newly inserted column
$ws.Range("B1").Value = "Week_Start_Date"

# Make sure the new column stores its values as plain text (so dates like
# "2025-01-05" are not auto-converted into Excel date serials).
$ws.Columns("B:B").NumberFormat = "@"

# 3. Week labels (column A) and week start dates (new column B)
$weeks = @(
    @{ Row = 2;  Label = "W1";  Date = "2025-01-05" },
    @{ Row = 3;  Label = "W2";  Date = "2025-01-12" },
    @{ Row = 4;  Label = "W3";  Date = "2025-01-19" },
    @{ Row = 5;  Label = "W4";  Date = "2025-01-26" },
    @{ Row = 6;  Label = "W5";  Date = "2025-02-02" },
    @{ Row = 7;  Label = "W6";  Date = "2025-02-09" },
    @{ Row = 8;  Label = "W7";  Date = "2025-02-16" },
    @{ Row = 9;  Label = "W8";  Date = "2025-02-23" },
    @{ Row = 10; Label = "W9";  Date = "2025-03-02" },
    @{ Row = 11; Label = "W10"; Date = "2025-03-09" },
    @{ Row = 12; Label = "W11"; Date = "2025-03-16" },
    @{ Row = 13; Label = "W12"; Date = "2025-03-23" },
    @{ Row = 14; Label = "W13"; Date = "2025-03-30" },
    @{ Row = 15; Label = "W14"; Date = "2025-04-06" },
    @{ Row = 16; Label = "W15"; Date = "2025-04-13" },
    @{ Row = 17; Label = "W16"; Date = "2025-04-20" }
)

foreach ($w in $weeks) {
    $r = $w.Row
    $ws.Range("A$r").Value = $w.Label
    $ws.Range("B$r").Value = $w.Date
    # is_holiday_week lives in column J after the insert; store as boolean FALSE
    $ws.Range("J$r").Value = $false
}
